# Pick-counts refresh: new pay-period data pulled in (adds several
# agents - e.g. ARJUNBHAI.PATEL, BOHD0676.KUSHLIAK, SURESH.DHAWAN - who
# picked up shifts late in the period) across every summary sheet in the
# workbook, per commit "added late to eff group".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PUTWALL PICKING")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "PutwallPickingQuantity"
$ws.Range("C1").Value = "UPH"
$data = New-Object "object[,]" 11,3
$data[0,0] = "ADOL798N.SEEMANNVAZQ"; $data[0,1] = 45.0; $data[0,2] = 13.98963730569948
$data[1,0] = "ARJUNBHAI.PATEL"; $data[1,1] = 27.0; $data[1,2] = 8.393782383419689
$data[2,0] = "BOHD0676.KUSHLIAK"; $data[2,1] = 161.0; $data[2,2] = 50.05181347150259
$data[3,0] = "DIAN4065.ENTRIALGO"; $data[3,1] = 27.0; $data[3,2] = 8.393782383419689
$data[4,0] = "LOANA.MBONGO"; $data[4,1] = 48.0; $data[4,2] = 14.92227979274611
$data[5,0] = "MICA0432.RIZKALLAMAR"; $data[5,1] = 52.0; $data[5,2] = 16.16580310880829
$data[6,0] = "PATR5027.AMEH"; $data[6,1] = 2.0; $data[6,2] = 0.6217616580310881
$data[7,0] = "SURESH.DHAWAN"; $data[7,1] = 99.0; $data[7,2] = 30.77720207253886
$data[8,0] = "THIE6554.DIALLO"; $data[8,1] = 106.0; $data[8,2] = 32.95336787564766
$data[9,0] = "TUSHAR.BHATIA"; $data[9,1] = 23.0; $data[9,2] = 7.150259067357513
$data[10,0] = "ZAHIDGUL.MINHAS"; $data[10,1] = 5.0; $data[10,2] = 1.55440414507772
$ws.Range("A2:C12").Value = $data

$ws = $wb.Worksheets.Item("REGULAR PICK")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "RegularPickQuantity"
$ws.Range("C1").Value = "UPH"
$data = New-Object "object[,]" 7,3
$data[0,0] = "ARJUNBHAI.PATEL"; $data[0,1] = 1.0; $data[0,2] = 0.310880829015544
$data[1,0] = "BOHD0676.KUSHLIAK"; $data[1,1] = 9.0; $data[1,2] = 2.797927461139896
$data[2,0] = "DIAN4065.ENTRIALGO"; $data[2,1] = 54.0; $data[2,2] = 16.78756476683938
$data[3,0] = "MARI882N.ABDELKADER"; $data[3,1] = 8.0; $data[3,2] = 2.487046632124352
$data[4,0] = "PATR5027.AMEH"; $data[4,1] = 4.0; $data[4,2] = 1.243523316062176
$data[5,0] = "WESL5337.CADETTE"; $data[5,1] = 30.0; $data[5,2] = 9.32642487046632
$data[6,0] = "ZAHIDGUL.MINHAS"; $data[6,1] = 4.0; $data[6,2] = 1.243523316062176
$ws.Range("A2:C8").Value = $data

$ws = $wb.Worksheets.Item("SINGLE PICK")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "SinglePickQuantity"
$ws.Range("C1").Value = "UPH"
$data = New-Object "object[,]" 7,3
$data[0,0] = "GIGNESH.PATEL"; $data[0,1] = 35.0; $data[0,2] = 10.88082901554404
$data[1,0] = "KADE3054.ZONGO"; $data[1,1] = 14.0; $data[1,2] = 4.352331606217616
$data[2,0] = "LOANA.MBONGO"; $data[2,1] = 50.0; $data[2,2] = 15.5440414507772
$data[3,0] = "SEPIDEH.AZARIHASHJIN"; $data[3,1] = 153.0; $data[3,2] = 47.56476683937824
$data[4,0] = "STAN9294.BAUER"; $data[4,1] = 86.0; $data[4,2] = 26.73575129533679
$data[5,0] = "TUSHAR.BHATIA"; $data[5,1] = 12.0; $data[5,2] = 3.730569948186528
$data[6,0] = "WESL5337.CADETTE"; $data[6,1] = 15.0; $data[6,2] = 4.66321243523316
$ws.Range("A2:C8").Value = $data

$ws = $wb.Worksheets.Item("REPLENISHMENT PICK")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "ReplenishmentPickQuantity"
$ws.Range("C1").Value = "UPH"
$data = New-Object "object[,]" 23,3
$data[0,0] = "AGNE8120.CARUTH"; $data[0,1] = 170.0; $data[0,2] = 52.84974093264248
$data[1,0] = "ARJUNBHAI.PATEL"; $data[1,1] = 82.0; $data[1,2] = 25.49222797927461
$data[2,0] = "ASHA1141.PAGE"; $data[2,1] = 23.0; $data[2,2] = 7.150259067357513
$data[3,0] = "BOHD0676.KUSHLIAK"; $data[3,1] = 112.0; $data[3,2] = 34.81865284974093
$data[4,0] = "DEVI789.SINGH"; $data[4,1] = 87.0; $data[4,2] = 27.04663212435233
$data[5,0] = "DIAN4065.ENTRIALGO"; $data[5,1] = 91.0; $data[5,2] = 28.29015544041451
$data[6,0] = "GIGNESH.PATEL"; $data[6,1] = 73.0; $data[6,2] = 22.69430051813471
$data[7,0] = "INUK4091.QAVAVAU"; $data[7,1] = 96.0; $data[7,2] = 29.84455958549223
$data[8,0] = "JEEW9554.SITUMUDALIG"; $data[8,1] = 84.0; $data[8,2] = 26.1139896373057
$data[9,0] = "KADE3054.ZONGO"; $data[9,1] = 50.0; $data[9,2] = 15.5440414507772
$data[10,0] = "LOANA.MBONGO"; $data[10,1] = 57.0; $data[10,2] = 17.72020725388601
$data[11,0] = "MARI882N.ABDELKADER"; $data[11,1] = 37.0; $data[11,2] = 11.50259067357513
$data[12,0] = "MICA0432.RIZKALLAMAR"; $data[12,1] = 103.0; $data[12,2] = 32.02072538860104
$data[13,0] = "NESR2403.ATTALAH"; $data[13,1] = 76.0; $data[13,2] = 23.62694300518135
$data[14,0] = "OMAR6689.KHAN"; $data[14,1] = 61.0; $data[14,2] = 18.96373056994819
$data[15,0] = "PRINCE.FORSON"; $data[15,1] = 101.0; $data[15,2] = 31.39896373056995
$data[16,0] = "STAN9294.BAUER"; $data[16,1] = 33.0; $data[16,2] = 10.25906735751295
$data[17,0] = "THIE6554.DIALLO"; $data[17,1] = 85.0; $data[17,2] = 26.42487046632124
$data[18,0] = "TUSHAR.BHATIA"; $data[18,1] = 99.0; $data[18,2] = 30.77720207253886
$data[19,0] = "WESL5337.CADETTE"; $data[19,1] = 83.0; $data[19,2] = 25.80310880829015
$data[20,0] = "WILDINE.JEUNE"; $data[20,1] = 172.0; $data[20,2] = 53.47150259067357
$data[21,0] = "YATI0689.YATIN"; $data[21,1] = 112.0; $data[21,2] = 34.81865284974093
$data[22,0] = "ZAKI0190.PHILLIPHORS"; $data[22,1] = 118.0; $data[22,2] = 36.68393782383419
$ws.Range("A2:C24").Value = $data

$ws = $wb.Worksheets.Item("QUICK MOVE")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "QuickMoveQuantity"
$ws.Range("C1").Value = "UPH"
$data = New-Object "object[,]" 13,3
$data[0,0] = "ADOL798N.SEEMANNVAZQ"; $data[0,1] = 173.0; $data[0,2] = 53.78238341968912
$data[1,0] = "BOHD0676.KUSHLIAK"; $data[1,1] = 79.0; $data[1,2] = 24.55958549222798
$data[2,0] = "DEVI789.SINGH"; $data[2,1] = 33.0; $data[2,2] = 10.25906735751295
$data[3,0] = "DIAN4065.ENTRIALGO"; $data[3,1] = 154.0; $data[3,2] = 47.87564766839378
$data[4,0] = "ESSE0616.UDEH"; $data[4,1] = 321.0; $data[4,2] = 99.79274611398964
$data[5,0] = "JEEW9554.SITUMUDALIG"; $data[5,1] = 10.0; $data[5,2] = 3.10880829015544
$data[6,0] = "MICA0432.RIZKALLAMAR"; $data[6,1] = 171.0; $data[6,2] = 53.16062176165803
$data[7,0] = "NESR2403.ATTALAH"; $data[7,1] = 97.0; $data[7,2] = 30.15544041450777
$data[8,0] = "STAN9294.BAUER"; $data[8,1] = 112.0; $data[8,2] = 34.81865284974093
$data[9,0] = "SURESH.DHAWAN"; $data[9,1] = 199.0; $data[9,2] = 61.86528497409326
$data[10,0] = "THIE6554.DIALLO"; $data[10,1] = 86.0; $data[10,2] = 26.73575129533679
$data[11,0] = "WESL5337.CADETTE"; $data[11,1] = 107.0; $data[11,2] = 33.26424870466321
$data[12,0] = "YATI0689.YATIN"; $data[12,1] = 42.0; $data[12,2] = 13.05699481865285
$ws.Range("A2:C14").Value = $data

$ws = $wb.Worksheets.Item("IDLE TIME")
$ws.Range("A1").Value = "UserID"
$ws.Range("B1").Value = "TotalIdleTime"
$data = New-Object "object[,]" 29,2
$data[0,0] = "ADOL798N.SEEMANNVAZQ"; $data[0,1] = 119.0
$data[1,0] = "AGNE8120.CARUTH"; $data[1,1] = 43.0
$data[2,0] = "ARJUNBHAI.PATEL"; $data[2,1] = 52.0
$data[3,0] = "BOHD0676.KUSHLIAK"; $data[3,1] = 50.0
$data[4,0] = "BUDD0680.TENNAKOON"; $data[4,1] = 47.0
$data[5,0] = "DEVI789.SINGH"; $data[5,1] = 52.0
$data[6,0] = "DIAN4065.ENTRIALGO"; $data[6,1] = 55.0
$data[7,0] = "ESSE0616.UDEH"; $data[7,1] = 109.0
$data[8,0] = "GIGNESH.PATEL"; $data[8,1] = 82.0
$data[9,0] = "INUK4091.QAVAVAU"; $data[9,1] = 120.0
$data[10,0] = "JEEW9554.SITUMUDALIG"; $data[10,1] = 28.0
$data[11,0] = "KADE3054.ZONGO"; $data[11,1] = 24.0
$data[12,0] = "LOANA.MBONGO"; $data[12,1] = 33.0
$data[13,0] = "MARI882N.ABDELKADER"; $data[13,1] = 131.0
$data[14,0] = "MICA0432.RIZKALLAMAR"; $data[14,1] = 26.0
$data[15,0] = "NESR2403.ATTALAH"; $data[15,1] = 60.0
$data[16,0] = "OMAR6689.KHAN"; $data[16,1] = 77.0
$data[17,0] = "PATR5027.AMEH"; $data[17,1] = 162.0
$data[18,0] = "PRINCE.FORSON"; $data[18,1] = 63.0
$data[19,0] = "SEPIDEH.AZARIHASHJIN"; $data[19,1] = 51.0
$data[20,0] = "STAN9294.BAUER"; $data[20,1] = 54.0
$data[21,0] = "SURESH.DHAWAN"; $data[21,1] = 133.0
$data[22,0] = "THIE6554.DIALLO"; $data[22,1] = 46.0
$data[23,0] = "TUSHAR.BHATIA"; $data[23,1] = 68.0
$data[24,0] = "WESL5337.CADETTE"; $data[24,1] = 47.0
$data[25,0] = "WILDINE.JEUNE"; $data[25,1] = 43.0
$data[26,0] = "YATI0689.YATIN"; $data[26,1] = 51.0
$data[27,0] = "ZAHIDGUL.MINHAS"; $data[27,1] = 140.0
$data[28,0] = "ZAKI0190.PHILLIPHORS"; $data[28,1] = 83.0
$ws.Range("A2:B30").Value = $data

$ws = $wb.Worksheets.Item("Total Units picked by hour")
$ws.Range("A1").Value = "Hour"
$ws.Range("B1").Value = "Regular Pick"
$ws.Range("C1").Value = "Single Pick"
$ws.Range("D1").Value = "Replenishment Pick"
$ws.Range("E1").Value = "Putwall Pick"
$data = New-Object "object[,]" 5,5
$data[0,0] = 20.0; $data[0,1] = -27.0; $data[0,2] = -24.0; $data[0,3] = -468.0; $data[0,4] = -1.0
$data[1,0] = 21.0; $data[1,1] = -41.0; $data[1,2] = -130.0; $data[1,3] = -788.0; $data[1,4] = -81.0
$data[2,0] = 22.0; $data[2,1] = -32.0; $data[2,2] = -168.0; $data[2,3] = -645.0; $data[2,4] = -357.0
$data[3,0] = 23.0; $data[3,1] = -10.0; $data[3,2] = -43.0; $data[3,3] = -104.0; $data[3,4] = -156.0
$data[4,0] = "Total"; $data[4,1] = -110.0; $data[4,2] = -365.0; $data[4,3] = -2005.0; $data[4,4] = -595.0
$ws.Range("A2:E6").Value = $data

